$p = $ppt.ActivePresentation
Write-Output "SlideMaster:"
try { Write-Output ($p.SlideMaster | Get-Member | Out-String) } catch { Write-Output "ERR: $_" }
